$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2023-06-19)
$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 71518.67920618469

# Row 3 (2023-03-21)
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 16.98373111632243
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 26.62400969366105
